# CSRinT3_Sp1_Annotation.xlsx bugfix:
#   - Measurement date field no longer carries a time component (it can be
#     supplied either in the "start time" field or in the raw data file).
#   - PO (Plant Ontology) growth-stage code is now stored as a proper
#     "PO:0007016" text id (was a bare, incorrectly-formatted number).
#   - Fixed a typo in the spectrometer channel name.
#   - Dropped the now-unused green "required" background from a few cells
#     that stopped needing manual data entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B4: measurement date -> drop the time-of-day portion ------------------
$ws.Range("B4").Value2 = "5/7/2012"

# --- D4: matching help text -> drop the trailing "time is in 24hr format" --
# D4 is a rich-text string (mixed italic run for "m/d/yyyy"); edit only the
# final (non-italic) run so the italic run in the middle survives intact.
$d4 = $ws.Range("D4")
$fullText = $d4.Text
$oldTail = ' with no leading zeros, e.g. "5/7/2012". The time is in 24hr format'
$newTail = ' with no leading zeros, e.g. "5/7/2012".'
$tailStart = $fullText.IndexOf($oldTail)
if ($tailStart -ge 0) {
    $tailChars = $d4.Characters($tailStart + 1, $oldTail.Length)
    $tailChars.Text = $newTail

    # Re-apply the italic run's own formatting (changing any run's .Text
    # collapses the cell to a single plain run, so restore the look of the
    # two surviving runs: the italic "m/d/yyyy" sample and the plain tail).
    $refreshed = $d4.Text
    $italicStart = $refreshed.IndexOf("m/d/yyyy")
    if ($italicStart -ge 0) {
        $italicChars = $d4.Characters($italicStart + 1, 8)
        $italicChars.Font.Name = "Calibri"
        $italicChars.Font.Size = 11
        $italicChars.Font.Color = 0
        $italicChars.Font.Italic = $true

        $tailPos = $italicStart + 8
        $tailLen = $refreshed.Length - $tailPos
        if ($tailLen -gt 0) {
            $tailChars2 = $d4.Characters($tailPos + 1, $tailLen)
            $tailChars2.Font.Name = "Calibri"
            $tailChars2.Font.Size = 11
            $tailChars2.Font.Color = 0
        }
    }
}

# --- B5: Growth Stage value -> real "PO:nnnnnnn" Plant Ontology id --------
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value2 = "PO:0007016"

# --- B6 / B7: start/end time -> no longer required, clear green fill -----
$ws.Range("B6").Interior.ColorIndex = -4142
$ws.Range("B7").Interior.ColorIndex = -4142

# --- B8: integration time -> no longer required, clear green fill --------
$ws.Range("B8").Interior.ColorIndex = -4142

# --- B10: fix "optimzed" -> "optimized" typo in channel name -------------
$ws.Range("B10").Value2 = "UCD_WUEoptimized_Channel1"

# --- refresh the on-screen selection to match the author's last edit -----
$ws.Range("D5").Select()
